$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) date serial value from 45183 (2023-09-14)
# to 45184 (2023-09-15) for rows 2-6, keeping the existing date formatting.
foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
